$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 4 columns: email_name | email_password | pop3svr | smtpsvr
# A new "imapsvr" field is being inserted as column C, pushing the existing
# pop3svr / smtpsvr columns one place to the right (-> D / E).

# 1) Shift the existing pop3svr/smtpsvr columns one column to the right.
$ws.Range("E1").Value = "smtpsvr"
$ws.Range("E2").Value = "smtp.rambler.ru"
$ws.Range("D1").Value = "pop3svr"
$ws.Range("D2").Value = "pop.rambler.ru"

# 2) Fill in the new imapsvr column in the freed-up column C.
$ws.Range("C1").Value = "imapsvr"
$ws.Range("C2").Value = "imap.rambler.ru"

# 3) Re-apply the existing header/data formatting to the new and shifted cells
#    (matching the look of the untouched B column in each row).
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)

# 4) Match the recorded selection left behind in the saved file.
$ws.Range("C6").Select()
